$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, new Timestamp (column A) value, new Load (column B) value
$data = @(
    @(2, 45454, 5780),
    @(3, 45454.01041666666, 5730),
    @(4, 45454.02083333334, 5680),
    @(5, 45454.03125, 5630),
    @(6, 45454.04166666666, 5590),
    @(7, 45454.05208333334, 5560),
    @(8, 45454.0625, 5530),
    @(9, 45454.07291666666, 5510),
    @(10, 45454.08333333334, 5500),
    @(11, 45454.09375, 5500),
    @(12, 45454.10416666666, 5500),
    @(13, 45454.11458333334, 5510),
    @(14, 45454.125, 5520),
    @(15, 45454.13541666666, 5520),
    @(16, 45454.14583333334, 5520),
    @(17, 45454.15625, 5520),
    @(18, 45454.16666666666, 5510),
    @(19, 45454.17708333334, 5510),
    @(20, 45454.1875, 5530),
    @(21, 45454.19791666666, 5590),
    @(22, 45454.20833333334, 5650),
    @(23, 45454.21875, 5740),
    @(24, 45454.22916666666, 5840),
    @(25, 45454.23958333334, 5950),
    @(26, 45454.25, 6130),
    @(27, 45454.26041666666, 6240),
    @(28, 45454.27083333334, 6340),
    @(29, 45454.28125, 6430),
    @(30, 45454.29166666666, 6500),
    @(31, 45454.30208333334, 6550),
    @(32, 45454.3125, 6570),
    @(33, 45454.32291666666, 6600),
    @(34, 45454.33333333334, 6610),
    @(35, 45454.34375, 6610),
    @(36, 45454.35416666666, 6600),
    @(37, 45454.36458333334, 6580),
    @(38, 45454.375, 6560),
    @(39, 45454.38541666666, 6540),
    @(40, 45454.39583333334, 6520),
    @(41, 45454.40625, 6510),
    @(42, 45454.41666666666, 6480),
    @(43, 45454.42708333334, 6470),
    @(44, 45454.4375, 6470),
    @(45, 45454.44791666666, 6470),
    @(46, 45454.45833333334, 6460),
    @(47, 45454.46875, 6450),
    @(48, 45454.47916666666, 6450),
    @(49, 45454.48958333334, 6450),
    @(50, 45454.5, 6460),
    @(51, 45454.51041666666, 6470),
    @(52, 45454.52083333334, 6470),
    @(53, 45454.53125, 6460),
    @(54, 45454.54166666666, 6440),
    @(55, 45454.55208333334, 6430),
    @(56, 45454.5625, 6430),
    @(57, 45454.57291666666, 6430),
    @(58, 45454.58333333334, 6450),
    @(59, 45454.59375, 6450),
    @(60, 45454.60416666666, 6450),
    @(61, 45454.61458333334, 6450),
    @(62, 45454.625, 6470),
    @(63, 45454.63541666666, 6480),
    @(64, 45454.64583333334, 6510),
    @(65, 45454.65625, 6530),
    @(66, 45454.66666666666, 6550),
    @(67, 45454.67708333334, 6590),
    @(68, 45454.6875, 6620),
    @(69, 45454.69791666666, 6670),
    @(70, 45454.70833333334, 6750),
    @(71, 45454.71875, 6820),
    @(72, 45454.72916666666, 6880),
    @(73, 45454.73958333334, 6950),
    @(74, 45454.75, 7010),
    @(75, 45454.76041666666, 7060),
    @(76, 45454.77083333334, 7100),
    @(77, 45454.78125, 7130),
    @(78, 45454.79166666666, 7170),
    @(79, 45454.80208333334, 7200),
    @(80, 45454.8125, 7240),
    @(81, 45454.82291666666, 7290),
    @(82, 45454.83333333334, 7330),
    @(83, 45454.84375, 7320),
    @(84, 45454.85416666666, 7310),
    @(85, 45454.86458333334, 7230),
    @(86, 45454.875, 7100),
    @(87, 45454.88541666666, 6980),
    @(88, 45454.89583333334, 6850),
    @(89, 45454.90625, 6690),
    @(90, 45454.91666666666, 6490),
    @(91, 45454.92708333334, 6340),
    @(92, 45454.9375, 6230),
    @(93, 45454.94791666666, 6120),
    @(94, 45454.95833333334, 5890),
    @(95, 45454.96875, 5850),
    @(96, 45454.97916666666, 5790),
    @(97, 45454.98958333334, 5730),
    @(98, 45455, 5680)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Cells.Item($row, 1).Value = $entry[1]
    $ws.Cells.Item($row, 2).Value = $entry[2]
}
